$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'63.211.16"
$ws.Cells.Item(2, 5).Value = "'  +0.79%  "
$ws.Cells.Item(3, 4).Value = "'3.023.83"
$ws.Cells.Item(3, 5).Value = "'  -2.45%  "
$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 5).Value = "'  -0.07%  "
$ws.Cells.Item(5, 4).Value = "'558.66"
$ws.Cells.Item(5, 5).Value = "'  +0.35%  "
$ws.Cells.Item(6, 4).Value = "'155.28"
$ws.Cells.Item(6, 5).Value = "'  -3.84%  "
$ws.Cells.Item(7, 4).Value = "'0.999"
$ws.Cells.Item(7, 5).Value = "'  -0.20%  "
$ws.Cells.Item(8, 4).Value = "'0.563"
$ws.Cells.Item(8, 5).Value = "'  -3.35%  "
$ws.Cells.Item(9, 4).Value = "'3.029.59"
$ws.Cells.Item(9, 5).Value = "'  -2.15%  "
$ws.Cells.Item(10, 5).Value = "'  -1.38%  "
$ws.Cells.Item(11, 4).Value = "'6.43"
$ws.Cells.Item(11, 5).Value = "'  -4.26%  "
$ws.Cells.Item(12, 4).Value = "'0.368"
$ws.Cells.Item(12, 5).Value = "'  -1.86%  "
$ws.Cells.Item(13, 4).Value = "'3.552.25"
$ws.Cells.Item(13, 5).Value = "'  -2.31%  "
$ws.Cells.Item(14, 5).Value = "'  -3.05%  "
$ws.Cells.Item(15, 4).Value = "'63.260.87"
$ws.Cells.Item(15, 5).Value = "'  +0.81%  "
$ws.Cells.Item(16, 4).Value = "'24.21"
$ws.Cells.Item(16, 5).Value = "'  -1.09%  "
$ws.Cells.Item(17, 2).Value = "'ShibaInu"
$ws.Cells.Item(17, 3).Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(17, 4).Value = "'0.0000151"
$ws.Cells.Item(17, 5).Value = "'  +0.29%  "
$ws.Cells.Item(18, 2).Value = "'WrappedEther"
$ws.Cells.Item(18, 3).Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(18, 4).Value = "'3.027.30"
$ws.Cells.Item(18, 5).Value = "'  -2.39%  "
$ws.Cells.Item(19, 4).Value = "'400.27"
$ws.Cells.Item(19, 5).Value = "'  +0.63%  "
$ws.Cells.Item(20, 4).Value = "'5.12"
$ws.Cells.Item(20, 5).Value = "'  +0.50%  "
$ws.Cells.Item(21, 4).Value = "'12.06"
$ws.Cells.Item(21, 5).Value = "'  -1.72%  "
$ws.Cells.Item(22, 4).Value = "'6.69"
$ws.Cells.Item(22, 5).Value = "'  -4.19%  "
$ws.Cells.Item(23, 4).Value = "'1.00"
$ws.Cells.Item(23, 5).Value = "'  +0.15%  "
$ws.Cells.Item(24, 4).Value = "'65.38"
$ws.Cells.Item(24, 5).Value = "'  -3.03%  "
$ws.Cells.Item(25, 2).Value = "'Polygon"
$ws.Cells.Item(25, 3).Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(25, 4).Value = "'0.467"
$ws.Cells.Item(25, 5).Value = "'  -1.25%  "
$ws.Cells.Item(26, 2).Value = "'Kaspa"
$ws.Cells.Item(26, 3).Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(26, 4).Value = "'0.189"
$ws.Cells.Item(26, 5).Value = "'  -3.61%  "
$ws.Cells.Item(27, 4).Value = "'0.0₃0992"
$ws.Cells.Item(27, 5).Value = "'  -0.67%  "
$ws.Cells.Item(28, 4).Value = "'8.73"
$ws.Cells.Item(28, 5).Value = "'  +1.78%  "
$ws.Cells.Item(29, 4).Value = "'0.997"
$ws.Cells.Item(29, 5).Value = "'  -0.54%  "
$ws.Cells.Item(30, 4).Value = "'1.00"
$ws.Cells.Item(30, 5).Value = "'  +0.07%  "
$ws.Cells.Item(31, 5).Value = "'  -0.17%  "
$ws.Cells.Item(32, 4).Value = "'20.42"
$ws.Cells.Item(32, 5).Value = "'  -1.50%  "
$ws.Cells.Item(33, 4).Value = "'161.15"
$ws.Cells.Item(33, 5).Value = "'  +5.57%  "
$ws.Cells.Item(34, 4).Value = "'4.78"
$ws.Cells.Item(34, 5).Value = "'  -0.01%  "
$ws.Cells.Item(35, 4).Value = "'1.12"
$ws.Cells.Item(35, 5).Value = "'  +2.99%  "
$ws.Cells.Item(36, 4).Value = "'6.08"
$ws.Cells.Item(36, 5).Value = "'  -0.94%  "
$ws.Cells.Item(37, 5).Value = "'  +0.89%  "
$ws.Cells.Item(38, 4).Value = "'2.549.78"
$ws.Cells.Item(38, 5).Value = "'  -5.41%  "
$ws.Cells.Item(39, 5).Value = "'  -2.99%  "
$ws.Cells.Item(40, 4).Value = "'22.97"
$ws.Cells.Item(40, 5).Value = "'  -0.78%  "
$ws.Cells.Item(41, 4).Value = "'3.96"
$ws.Cells.Item(41, 5).Value = "'  -0.45%  "
$ws.Cells.Item(42, 4).Value = "'37.75"
$ws.Cells.Item(42, 5).Value = "'  -0.93%  "
$ws.Cells.Item(43, 4).Value = "'0.671"
$ws.Cells.Item(43, 5).Value = "'  -2.44%  "
$ws.Cells.Item(44, 4).Value = "'0.0603"
$ws.Cells.Item(44, 5).Value = "'  +0.93%  "
$ws.Cells.Item(45, 4).Value = "'0.0251"
$ws.Cells.Item(45, 5).Value = "'  -0.42%  "
$ws.Cells.Item(46, 4).Value = "'5.10"
$ws.Cells.Item(46, 5).Value = "'  -1.66%  "
$ws.Cells.Item(47, 4).Value = "'0.997"
$ws.Cells.Item(47, 5).Value = "'  -0.30%  "
$ws.Cells.Item(48, 4).Value = "'20.35"
$ws.Cells.Item(48, 5).Value = "'  -0.94%  "
$ws.Cells.Item(49, 4).Value = "'269.73"
$ws.Cells.Item(49, 5).Value = "'  -3.94%  "
$ws.Cells.Item(50, 2).Value = "'Stellar"
$ws.Cells.Item(50, 3).Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(50, 4).Value = "'0.0948"
$ws.Cells.Item(50, 5).Value = "'  -2.11%  "
$ws.Cells.Item(51, 2).Value = "'WhiteBITCoin"
$ws.Cells.Item(51, 3).Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Cells.Item(51, 4).Value = "'10.47"
$ws.Cells.Item(51, 5).Value = "'  +0.18%  "
